$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'46.637.12"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "'  +6.31%  "
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = "'2.294.08"
$ws.Range('D3').ClearFormats()
$ws.Range('E4').Value = "'  +0.08%  "
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = "'305.29"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "'  +2.04%  "
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = "'100.29"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = "'  +11.25%  "
$ws.Range('E6').ClearFormats()
$ws.Range('D7').Value = "'0.566"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "'  +2.20%  "
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = "'1.00"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = "'  +0.02%  "
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = "'0.518"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "'  +5.72%  "
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = "'36.60"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = "'  +11.48%  "
$ws.Range('E10').ClearFormats()
$ws.Range('D11').Value = "'0.0789"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = "'  +1.78%  "
$ws.Range('E11').ClearFormats()
$ws.Range('E12').Value = "'  +6.17%  "
$ws.Range('E12').ClearFormats()
$ws.Range('E13').Value = "'  +0.37%  "
$ws.Range('E13').ClearFormats()
$ws.Range('D14').Value = "'2.644.37"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "'  +3.45%  "
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = "'2.290.08"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "'  +3.26%  "
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = "'13.80"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = "'  +3.27%  "
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = "'0.809"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "'  +4.52%  "
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = "'46.630.95"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "'  +6.65%  "
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = "'12.98"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "'  +13.61%  "
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = "'0.0₃0936"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "'  +3.78%  "
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = "'6.01"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = "'  +1.24%  "
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = "'66.33"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "'  +3.01%  "
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = "'247.45"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "'  +4.89%  "
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = "'2.90"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = "'  +3.36%  "
$ws.Range('E24').ClearFormats()
$ws.Range('E25').Value = "'  -0.01%  "
$ws.Range('E25').ClearFormats()
$ws.Range('E26').Value = "'  +3.10%  "
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = "'42.93"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "'  +12.12%  "
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = "'2.23"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "'  +1.10%  "
$ws.Range('E28').ClearFormats()
$ws.Range('D29').Value = "'9.84"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "'  +5.15%  "
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = "'19.90"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "'  +3.87%  "
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = "'2.82"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = "'  +13.60%  "
$ws.Range('E31').ClearFormats()
$ws.Range('E32').Value = "'  +4.59%  "
$ws.Range('E32').ClearFormats()
$ws.Range('D33').Value = "'147.10"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "'  -4.02%  "
$ws.Range('E33').ClearFormats()
$ws.Range('D34').Value = "'0.0792"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = "'  +4.59%  "
$ws.Range('E34').ClearFormats()
$ws.Range('D36').Value = "'0.115"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "'  +11.69%  "
$ws.Range('E36').ClearFormats()
$ws.Range('D37').Value = "'0.117"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "'  +0.68%  "
$ws.Range('E37').ClearFormats()
$ws.Range('B38').Value = "'ARBITRUM"
$ws.Range('B38').ClearFormats()
$ws.Range('C38').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('C38').ClearFormats()
$ws.Range('D38').Value = "'1.77"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = "'  +5.72%  "
$ws.Range('E38').ClearFormats()
$ws.Range('B39').Value = "'Celestia"
$ws.Range('B39').ClearFormats()
$ws.Range('C39').Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range('C39').ClearFormats()
$ws.Range('D39').Value = "'15.86"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "'  +20.92%  "
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = "'4.04"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "'  +11.87%  "
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = "'3.36"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "'  +6.08%  "
$ws.Range('E41').ClearFormats()
$ws.Range('E42').Value = "'  +1.58%  "
$ws.Range('E42').ClearFormats()
$ws.Range('E43').Value = "'  +0.03%  "
$ws.Range('E43').ClearFormats()
$ws.Range('D44').Value = "'1.95"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "'  +10.36%  "
$ws.Range('E44').ClearFormats()
$ws.Range('D45').Value = "'1.823.80"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "'  -0.56%  "
$ws.Range('E45').ClearFormats()
$ws.Range('D46').Value = "'86.85"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "'  +19.59%  "
$ws.Range('E46').ClearFormats()
$ws.Range('D47').Value = "'0.194"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "'  +6.93%  "
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = "'72.76"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = "'  +7.44%  "
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = "'4.89"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = "'  +7.09%  "
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = "'95.49"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "'  +1.03%  "
$ws.Range('E50').ClearFormats()
$ws.Range('D51').Value = "'2.521.32"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "'  +3.38%  "
$ws.Range('E51').ClearFormats()
